# data cleanup continued in player_per_game_df
# Insert two new player rows (LaMelo Ball, LeBron James) into the sorted
# player/award pivot table, each with an award count of 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "LaMelo Ball" row just above the existing "Larry Bird" row (row 38),
# pushing Larry Bird and everything after it down by one row.
$ws.Rows.Item(38).Insert()
$ws.Range("A38").Value = "LaMelo Ball"
$ws.Range("B38").Value = 2

# Insert "LeBron James" row just above the existing "Luka Dončić" row.
# After the first insert, Luka Dončić now sits at row 41, so insert above it.
$ws.Rows.Item(41).Insert()
$ws.Range("A41").Value = "LeBron James"
$ws.Range("B41").Value = 2
